# Update "想去人数" (interested-count) figures on both the "展览" and
# "全部类型" worksheets to reflect the latest generated output.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F8").Value = 66
    $ws.Range("F10").Value = 5633
    $ws.Range("F11").Value = 5001
}
